$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated cell values (Coin price / 1h volume change columns).
# Some "Price" values (e.g. "1.000", "0.7269") read as plain numbers to Excel's
# auto-detection, so for those cells we briefly force a Text number format,
# write the literal string, then restore the default "Normal" style so the
# cell formatting matches the source workbook exactly.

$ws.Range('D2').Value = '29.308.16'
$ws.Range('D3').Value = '1.868.10'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7269'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.43%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '241.07'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.26%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07890'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.60%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3094'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.41%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '25.28'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.14%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08247'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').Value = '1.876.94'
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.7236'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.38%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.244'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.48%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '90.81'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.76%  '
$ws.Range('D16').Value = '29.368.71'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.866'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.99%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '244.12'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.16%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007835'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.24'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.39%  '
$ws.Range('D21').Value = '2.120.42'
$ws.Range('E21').Value = '  -0.47%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.0000'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.023'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +6.69%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.9996'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1601'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +12.72%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '162.30'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.959'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.27'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.11%  '
$ws.Range('E30').Value = '  +1.66%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.405'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.49%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.113'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.90%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05215'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.92%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.938'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.80%  '
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7291'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +3.41%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.678'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  +0.89%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.701'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('D40').Value = '1.164.54'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9062'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.39%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.105'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.32%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '72.50'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.9997'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '101.93'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').Value = '2.018.67'
$ws.Range('E46').Value = '  -0.71%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5279'
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.778'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.46%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.891'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +5.85%  '
$ws.Range('E50').Value = '  +1.54%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4276'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.70%  '
